$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Stats")

# End-of-season round up: clear out the last mini "player of the week" style
# table (rows 353-377, columns B-F) that listed one-off leftover entries,
# while keeping the existing cell formatting intact.
$ws.Range("B353:F377").ClearContents()

# Update the selection to reflect the (slightly) extended range left selected
# after the clean-up.
$ws.Range("B353:F379").Select()
